# Regenerate quadratic/linear problem values
# (commit: "volver a generar problemas cuadraticos y lineales")
#
# Values that look like plain numbers must be entered with a leading
# apostrophe so Excel keeps storing them as text (shared-string) cells,
# exactly as they were in the original workbook, instead of silently
# converting them to numeric cells.
#
# NOTE: the workbook has two worksheets whose names differ only by case
# ("Vector_bf" and "Vector_BF"). Worksheet name lookups are case
# insensitive, so both names would resolve to the same sheet. Numeric
# sheet indices are used instead to unambiguously address each one:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = "1.7500000000000169 - 2x_1 + y_1 - y_2"
$ws.Range("B2").Value = "'0.7499999999999831"
$ws.Range("D2").Value = "'0.77"
$ws.Range("E2").Value = "'3.5"
$ws.Range("F2").Value = "'1.1"

$ws.Range("A3").Value = "6.549999999999994 + x_1 - 3x_2 + y_2"
$ws.Range("B3").Value = "'-8.549999999999994"
$ws.Range("D3").Value = "'0.46"
$ws.Range("E3").Value = "'9.9"
$ws.Range("F3").Value = "'4.3"

$ws.Range("A4").Value = "104.6 - y_1"
$ws.Range("B4").Value = "'-104.6"
$ws.Range("D4").Value = "'0.41"
$ws.Range("E4").Value = "'0.6"
$ws.Range("F4").Value = "'8.0"

$ws.Range("A5").Value = "-2.05 - y_2"
$ws.Range("B5").Value = "'-2.05"
$ws.Range("D5").Value = "'0.64"
$ws.Range("E5").Value = "'4.5"
$ws.Range("F5").Value = "'8.100000000000001"

# --- Punto_modificado ---------------------------------------------------
$ws = $wb.Worksheets.Item(4)

$ws.Range("A2").Value = "'52.150000000000006"
$ws.Range("B2").Value = "'20.25"
$ws.Range("C2").Value = "'104.6"
$ws.Range("D2").Value = "'2.05"

# --- Vector_bf ------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)

$ws.Range("A2").Value = "'3.64"
$ws.Range("A3").Value = "'-0.050000000000000044"

# --- Vector_BF ------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)

$ws.Range("A2").Value = "'-0.9000000000000004"
$ws.Range("A3").Value = "'28.700000000000003"
$ws.Range("A4").Value = "'-3.4"
$ws.Range("A5").Value = "'-1.9000000000000004"
